$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename task in row 5 (previously "setup environtment backend(routing dll)")
$ws.Range("C5").Value = "setup backend(routing dll)"

# Update status of "web interface" row (row 6) from ON PROGRESS to DONE
$ws.Range("F6").Value = "DONE"

# Remove the borders that were applied to the C3:E9 block (task/time columns)
$ws.Range("C3:E9").Borders.LineStyle = -4142

# Restore default view state (no frozen/scrolled top-left cell, selection on E26)
$ws.Range("E26").Select()

$wb.Save()
